$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'87.992.87"
$ws.Range("E2").Value = "  -2.02%  "

$ws.Range("D3").Value = "'3.060.11"
$ws.Range("E3").Value = "  -4.55%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").Value = "'209.91"
$ws.Range("E5").Value = "  -3.78%  "

$ws.Range("D6").Value = "'618.95"
$ws.Range("E6").Value = "  -1.64%  "

$ws.Range("D7").Value = "'0.371"
$ws.Range("E7").Value = "  -6.05%  "

$ws.Range("D8").Value = "'0.798"
$ws.Range("E8").Value = "  +14.42%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").Value = "'3.056.41"
$ws.Range("E10").Value = "  -4.67%  "

$ws.Range("E11").Value = "  +2.74%  "

$ws.Range("D12").Value = "'0.179"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("D13").Value = "'0.0000238"
$ws.Range("E13").Value = "  -9.67%  "

$ws.Range("D14").Value = "'5.28"
$ws.Range("E14").Value = "  -2.83%  "

$ws.Range("D15").Value = "'87.759.11"
$ws.Range("E15").Value = "  -2.15%  "

$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'31.86"
$ws.Range("E16").Value = "  -5.99%  "

$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "'3.617.66"
$ws.Range("E17").Value = "  -4.84%  "

$ws.Range("D18").Value = "'3.048.70"
$ws.Range("E18").Value = "  -4.98%  "

$ws.Range("D19").Value = "'3.25"
$ws.Range("E19").Value = "  -5.54%  "

$ws.Range("D20").Value = "'0.0000200"
$ws.Range("E20").Value = "  -13.31%  "

$ws.Range("D21").Value = "'13.22"
$ws.Range("E21").Value = "  -2.97%  "

$ws.Range("D22").Value = "'419.74"
$ws.Range("E22").Value = "  -4.61%  "

$ws.Range("D23").Value = "'8.14"
$ws.Range("E23").Value = "  -6.53%  "

$ws.Range("D24").Value = "'4.88"
$ws.Range("E24").Value = "  -5.06%  "

$ws.Range("D25").Value = "'5.48"
$ws.Range("E25").Value = "  +2.38%  "

$ws.Range("D26").Value = "'11.73"
$ws.Range("E26").Value = "  -3.52%  "

$ws.Range("D27").Value = "'81.79"
$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("E28").Value = "  -6.09%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("E30").Value = "  +8.59%  "

$ws.Range("D31").Value = "'0.170"
$ws.Range("E31").Value = "  +5.75%  "

$ws.Range("D32").Value = "'8.03"
$ws.Range("E32").Value = "  -6.30%  "

$ws.Range("D33").Value = "'507.42"
$ws.Range("E33").Value = "  -7.81%  "

$ws.Range("D34").Value = "'3.57"
$ws.Range("E34").Value = "  -14.20%  "

$ws.Range("D35").Value = "'6.74"
$ws.Range("E35").Value = "  -5.57%  "

$ws.Range("D36").Value = "'1.79"
$ws.Range("E36").Value = "  -7.15%  "

$ws.Range("D37").Value = "'1.23"
$ws.Range("E37").Value = "  -7.47%  "

$ws.Range("D38").Value = "'22.21"
$ws.Range("E38").Value = "  -1.54%  "

$ws.Range("D39").Value = "'0.130"
$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").Value = "'22.19"
$ws.Range("E40").Value = "  -0.91%  "

$ws.Range("E41").Value = "  +0.12%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "'0.360"
$ws.Range("E43").Value = "  -4.94%  "

$ws.Range("D44").Value = "'147.28"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.80"
$ws.Range("E45").Value = "  -7.86%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.132"
$ws.Range("E46").Value = "  +5.47%  "

$ws.Range("D47").Value = "'43.32"
$ws.Range("E47").Value = "  -1.31%  "

$ws.Range("D48").Value = "'0.0677"
$ws.Range("E48").Value = "  +11.30%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'155.70"
$ws.Range("E49").Value = "  -10.96%  "

$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.18"
$ws.Range("E50").Value = "  -6.46%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.699"
$ws.Range("E51").Value = "  -10.03%  "
